# Logic tree input file updated
# Insert two new "Possible_Problem" answer rows into the decision tree:
#   - one right after the "Does the engine run too hot?" question block starts (new row 13)
#   - one right after the "Are you seeing the coolant on the ground..." question block starts (new row 18)
# Every other existing row shifts down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionHot = "Problem:Does the engine run too hot ? (check temperature gauge on the dashboard) (Please answer as: Yes, No, Not Sure, Fluctuates)"
$questionSeeing = "Problem:Are you seeing the coolant on the ground every time you park? (Please answer as: Yes, No)"
$possibleProblemLabel = "Possible_Problem"
$possibleProblemDetail = "Possible_Problem:50% Cylinder Head Gaskets`n20% Radiator (Hoses)`n20% Water Pump`n10% Thermostat housing"

# Insert a new row before current row 13 (shifts old rows 13-21 down by one)
$ws.Rows.Item(13).Insert()
$ws.Cells.Item(13, 1).Value = $questionHot
$ws.Cells.Item(13, 2).Value = $possibleProblemLabel
$ws.Cells.Item(13, 3).Value = $possibleProblemDetail
$ws.Cells.Item(13, 3).WrapText = $true
$ws.Rows.Item(13).RowHeight = 230.4

# Insert a new row before current row 18 (which, after the previous insert, is the first
# row of the "Are you seeing the coolant..." block) - shifts those rows down by one more
$ws.Rows.Item(18).Insert()
$ws.Cells.Item(18, 1).Value = $questionSeeing
$ws.Cells.Item(18, 2).Value = $possibleProblemLabel
$ws.Cells.Item(18, 3).Value = $possibleProblemDetail
$ws.Cells.Item(18, 3).WrapText = $true
$ws.Rows.Item(18).RowHeight = 230.4

# Update view state to match final selection in the edited file
$ws.Application.ActiveWindow.ScrollRow = 23
$ws.Range("B23").Select()
